$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 18 ---
$ws.Range("A17").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A18").Value = 43891
$ws.Range("B18").Value = "1st review"
$ws.Range("C18").Value = "1st review"

# --- Row 19 ---
$ws.Range("A17").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A19").Value = 43983
$ws.Range("B19").Value = "UI DESIGNE START"
$ws.Range("C19").Value = "UI DESIGN START"

# --- Row 20 ---
$ws.Range("A17").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A20").Value = 44013
$ws.Range("B20").Value = "UI design was completed"
$ws.Range("C20").Value = "ui designed was completed."

# --- back-fill D19, D20 ---
$ws.Range("D19").Value = "java class was taken of different 14 topics"
$ws.Range("D20").Value = "spring  boot class was taken."

# --- jump to row 23, the date-looking text in A23 (plain text, not a real date) ---
$ws.Range("A23").Value = "10/01/202"

# --- back to D21 ---
$ws.Range("D21").Value = "try to implement  the ui design  in eclipse."

# --- finish remaining cells: row 21, row 22, and B23 ---
$ws.Range("A17").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A21").Value = 44044
$ws.Range("B21").Value = "java class was taken of different 14 topics"

$ws.Range("A17").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A22").Value = 44075
$ws.Range("B22").Value = "spring  boot class was taken."

$ws.Range("B23").Value = "try to implement  the ui design  in eclipse."

# --- update the view: select A17, clear any left-scrolled topLeftCell ---
$ws.Range("A17").Select()
